$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.762.32'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.147.04'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '578.94'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').Value = '149.04'
$ws.Range('E6').Value = '  -1.16%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.144.44'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('E10').Value = '  -2.54%  '
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').Value = '0.0000263'
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('D14').Value = '37.12'
$ws.Range('E14').Value = '  -2.64%  '
$ws.Range('D15').Value = '3.663.07'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').Value = '64.779.50'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '3.162.42'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').Value = '7.14'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').Value = '504.56'
$ws.Range('E20').Value = '  -2.48%  '
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('E22').Value = '  -2.91%  '
$ws.Range('D23').Value = '15.16'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('D25').Value = '84.38'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('E28').Value = '  +1.97%  '
$ws.Range('D29').Value = '2.18'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '27.58'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').Value = '2.78'
$ws.Range('E31').Value = '  +3.46%  '
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('D33').Value = '1.20'
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('D34').Value = '6.24'
$ws.Range('E34').Value = '  +2.06%  '
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('D36').Value = '54.87'
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('D37').Value = '484.79'
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('E38').Value = '  +2.70%  '
$ws.Range('D39').Value = '0.0415'
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  -1.98%  '
$ws.Range('D41').Value = '8.75'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').Value = '2.990.04'
$ws.Range('E42').Value = '  -4.06%  '
$ws.Range('D43').Value = '0.117'
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').Value = '0.281'
$ws.Range('E45').Value = '  -4.39%  '
$ws.Range('D46').Value = '28.08'
$ws.Range('E46').Value = '  -4.03%  '
$ws.Range('D47').Value = '0.0₃0586'
$ws.Range('E47').Value = '  +1.42%  '
$ws.Range('D49').Value = '0.113'
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('E50').Value = '  -3.14%  '
$ws.Range('D51').Value = '2.48'
$ws.Range('E51').Value = '  +14.08%  '
